$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: new record "Rastgele birisi"
$ws.Range("A8").Value = 2313213123
$ws.Range("A8").Locked = $true
$ws.Range("B8").Value = "Rastgele birisi"
$ws.Range("C8").Value = "Yazılımcı,Yardımcı"

# Row 9: new record "Birisi daha"
$ws.Range("A9").Value = 6546365
$ws.Range("B9").Value = "Birisi daha"
$ws.Range("C9").Value = "Takım üyesi"
